$d = $word.ActiveDocument

# Target the first paragraph, which holds the **ID__...__ID** marker text.
$para = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with a 5pt (space-only) margin,
# matching the new <w:pBdr> with w:space="5" on each edge.
$para.Borders.DistanceFromTop = 5
$para.Borders.DistanceFromLeft = 5
$para.Borders.DistanceFromBottom = 5
$para.Borders.DistanceFromRight = 5

# Bump the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$para.LeftIndent = 11.25

# Replace the marker text (which previously spanned two runs - the ID text
# plus a trailing space run) with the updated, single-run marker text.
$para.Range.Find.Execute("**ID__AFFARS_mp_5315_3_topic_33__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP_5315_3_6_3__ID**", 2)
